$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old rows 2-5 (unit2 dictation rows), row 1 stays as-is
$ws.Range("A2:J5").EntireRow.Delete() | Out-Null

# New K1:N1 placeholder cells ("-")
$ws.Range("K1").Value = "-"
$ws.Range("L1").Value = "-"
$ws.Range("M1").Value = "-"
$ws.Range("N1").Value = "-"

# Row 2
$ws.Range("A2").Value = "A1"
$ws.Range("B2").Value = "Reading "
$ws.Range("C2").Value = "Unit1"
$ws.Range("D2").Value = "In this lesson you will read a small text . When you are ready to start press the button and speak into the microphone . Keep trying until you get it right ! Don’t worry about making mistakes , mistakes mean you are improving ​"
$ws.Range("E2").Value = "ستسمع في هذا الدرس بعض الكلمات، حاول كتابتها بشكل صحيح. يمكنك سماع كل كلمة كم مرة تريد. تم أخذ كل هذه الكلمات لذا لن تجدها صعبة للغاية (إذا كنت قد قمت بتدوين الملاحظات:)​"
$ws.Range("F2").Value = "-"
$ws.Range("G2").Value = "من النص هل يمكنك ايجاد الكلمة الانجليزية للكلامات الاتية ؟ ​"
$ws.Range("H2").Value = "صداقة,الصبغ او الرسم بالالوان,فرشاة او ريشة,عبر او خلال,حقل,نزهة,لذيذة,بطيخة,عنب,قوس قزح,يقول,يلعب "
$ws.Range("I2").Value = "Friendship,Paint,Brush,Through,Field,Picnic,Yummy,Watermelon,Grapes,Rainbow,Says,Play"
$ws.Range("J2").Value = "findWordsFromPassage"
$ws.Range("K2").Value = "Friendship"
$ws.Range("L2").Value = "الصداقة "
$ws.Range("M2").Value = "In a little village, there were three friends: Red, Blue, and Green.​  Red likes to play with his red ball. Blue likes to paint with her blue brush. And Green likes to run through the green fields.​  They want to have a picnic. Red brings juicy red apples, Blue brings yummy watermelon and Green brings crunchy green grapes.​  They sit on a colorful blanket and enjoy their picnic under the bright sun. Suddenly, they see a rainbow in the sky.​  ""It's so pretty!"" says Blue.​  ""I see red, blue, and green!"" says Red.​  ""And yellow, purple, and orange too!"" adds Green.​  They laugh and play, feeling happy and colorful together. And from that day on, they know that friendship was the brightest color of all."
$ws.Range("N2").Value = "في قرية صغيرة، كان هناك ثلاثة أصدقاء: الأحمر والأزرق والأخضر.​  يحب الأحمر اللعب بكرته الحمراء. يحب الأزرق الرسم بفرشاته الزرقاء. والأخضر يحب الركض عبر الحقول الخضراء.​  يريدون الذهاب الى نزهة . الأحمر يجلب التفاح الأحمر العصير، والأزرق يجلب البطيخ اللذيذ والأخضر يجلب العنب الأخضر المقرمش. ​  يجلسون على بطانية ملونة ويستمتعون بنزهة تحت أشعة الشمس الساطعة. وفجأة رأوا قوس قزح في السماء.​  ""انها جميلة جدا!"" يقول الأزرق.​  ""أرى الأحمر والأزرق والأخضر!"" يقول الأحمر.​  ""والأصفر والأرجواني والبرتقالي أيضًا!"" يضيف الأخضر.​  إنهم يضحكون ويلعبون، ويشعرون بالسعادة والألوان معًا.​   ومن ذلك اليوم فصاعدًا، عرفوا أن الصداقة كانت ألمع الألوان على الإطلاق."

# Row 3
$ws.Range("A3").Value = "A1"
$ws.Range("B3").Value = "Reading "
$ws.Range("C3").Value = "Unit1"
$ws.Range("D3").Value = "In this lesson you will read a small text . When you are ready to start press the button and speak into the microphone . Keep trying until you get it right ! Don’t worry about making mistakes , mistakes mean you are improving ​"
$ws.Range("E3").Value = "ستسمع في هذا الدرس بعض الكلمات، حاول كتابتها بشكل صحيح. يمكنك سماع كل كلمة كم مرة تريد. تم أخذ كل هذه الكلمات لذا لن تجدها صعبة للغاية (إذا كنت قد قمت بتدوين الملاحظات:)​"
$ws.Range("F3").Value = "-"
$ws.Range("G3").Value = "جاوب على الاسئلة الاتية و قارن بنفسك الاجوبة . لا تنسى كتابة الاجوبة الصحيحة على الدفتر ! حاول بقدر ما تسطيع الاجابة , نعلم انها صعبة قليلا لكن بالصعوبة نترقى ! ​"
$ws.Range("H3").Value = "What is the name of the speaker? ما هو اسم المتحدث؟​,How old is the speaker? كم عمر المتحدث؟​,How many birds were there in the sky ? كم كان عدد الطيور في السماء؟​,How many stuffed animals does the speaker have? كم عدد الحيوانات المحشية التي يملكها المتحدث؟,What does the speaker like to do ? ماذا يحب المتحدث أن يفعل؟ ​"
$ws.Range("I3").Value = "The name of the speaker is Emily,The speaker is six years old ​,There were 5 birds flying in the sky​,The speaker has 5 stuffed animals​,She loves spending time with my toys and creating new adventures for them to go on!"
$ws.Range("J3").Value = "answerQuestionsFromPassage"
$ws.Range("K3").Value = "Emily "
$ws.Range("L3").Value = "إيميلي ​"
$ws.Range("M3").Value = "Hello! My name is Emily. I am six years old. I have three red balloons. I like red because it's my favorite color. I also have two blue toy cars and one green teddy bear. My teddy bear's name is Teddy. He's my best friend. I like to play with my toys in the park. Yesterday, I saw five birds flying in the sky. They were colorful: two were blue, two were yellow, and one was green. It was a beautiful sight! There were also three colorful balloons floating in the air: one red, one blue, and one yellow. After painting, I played with my favorite toys. I have five stuffed animals: a brown bear, a white bunny, a black and white panda, a pink elephant, and a purple unicorn. Each of them has their own special spot on my bed. I love spending time with my toys and creating new adventures for them to go on!"
$ws.Range("N3").Value = "مرحبًا! إسمي إيميلي. أنا ست سنوات من العمر. لدي ثلاث بالونات حمراء. أحب اللون الأحمر لأنه لوني المفضل. لدي أيضًا سيارتان لعبة باللون الأزرق ودبدوب أخضر. اسم الدبدوب الخاص بي هو تيدي. انه صديقي المفضل. أحب اللعب بألعابي في الحديقة. بالأمس رأيت خمسة طيور تحلق في السماء. كانت ملونة: اثنان باللون الأزرق، واثنان باللون الأصفر، وواحد باللون الأخضر. كان منظرا جميلا! كانت هناك أيضًا ثلاثة بالونات ملونة تطفو في الهواء: واحدة حمراء، وواحدة زرقاء، وواحدة صفراء. بعد الرسم، لعبت بألعابي المفضلة. لدي خمسة حيوانات محشوة: دب بني، وأرنب أبيض، وباندا أبيض وأسود، وفيل وردي، ووحيد قرن أرجواني. كل واحد منهم لديه مكانه الخاص على سريري. أحب قضاء الوقت مع ألعابي وخلق مغامرات جديدة لهم للاستمرار فيها!"

# Column M width
$ws.Columns.Item(13).ColumnWidth = 16.1640625

# View state
$null = $ws.Range("H8").Select()
$excel.ActiveWindow.Zoom = 137

